$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 2 data values to be numeric (previously stored as shared strings)
$ws.Range("A2").Value = 1000
$ws.Range("B2").Value = 20
$ws.Range("C2").Value = 0.5
$ws.Range("D2").Value = 0.5
$ws.Range("E2").Value = 0.1
$ws.Range("F2").Value = 10
$ws.Range("G2").Value = 10
$ws.Range("H2").Value = 1
